$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of the data block (row 353),
# pushing the existing rows 353-378 down by one (to 354-379).
$ws.Rows.Item(353).Insert()

$ws.Range("A353").Value = 4
$ws.Range("B353").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C353").Value = "Los Lagos"
$ws.Range("D353").Value = 45021
$ws.Range("E353").Value = 10
$ws.Range("F353").Value = 100112032
$ws.Range("G353").Value = "Zapallo italiano"
$ws.Range("H353").Value = "Sin especificar"
$ws.Range("I353").Value = "Primera"
$ws.Range("J353").Value = 35
$ws.Range("K353").Value = 11000
$ws.Range("L353").Value = 11000
$ws.Range("M353").Value = 11000
$ws.Range("N353").Value = "$/caja 50 unidades"
$ws.Range("O353").Value = "Región Metropolitana"
$ws.Range("P353").Value = 220
$ws.Range("Q353").Value = 50
$ws.Range("R353").Value = "Hortaliza"
